# Rows 2-11 of the "Artfynd" sheet are rotated: the block of the last four
# records (old rows 8-11) moves to the top (new rows 2-5), and the first six
# records (old rows 2-7) shift down underneath them (new rows 6-11).
# Row-by-row this is the mapping new-row <- old-row:
#   2<-8  3<-9  4<-10  5<-11  6<-2  7<-3  8<-4  9<-5  10<-6  11<-7
# Rows 12-14 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"

# 1) Snapshot every source row (2..11) BEFORE any writes happen, since the
#    source and destination rows overlap (this is a rotation in place).
$rowData = @{}
for ($r = 2; $r -le 11; $r++) {
    $rng = $ws.Range("$firstCol$r`:$lastCol$r")
    $rowData[$r] = $rng.Value2
}

# 2) New-row <- old-row mapping describing the rotation.
$mapping = @{
    2  = 8
    3  = 9
    4  = 10
    5  = 11
    6  = 2
    7  = 3
    8  = 4
    9  = 5
    10 = 6
    11 = 7
}

# 3) Columns Y (Startdatum) and AA (Slutdatum) hold plain text that looks
#    like dates ("2020-09-26"); force those destination cells to Text format
#    first so Value2 doesn't silently reinterpret them as date serials.
for ($destRow = 2; $destRow -le 11; $destRow++) {
    $ws.Range("Y$destRow").NumberFormat = "@"
    $ws.Range("AA$destRow").NumberFormat = "@"
}

# 4) Write each source row's snapshot into its new destination row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $destRng = $ws.Range("$firstCol$destRow`:$lastCol$destRow")
    $destRng.Value2 = $rowData[$srcRow]
}
